# Update the "Förändrad" (Changed) date column (C) for all data rows.
# Every value in column C (rows 2-275) was 46060 (2026-02-07) and is
# bumped forward by one day to 46061 (2026-02-08).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 275 }

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    $current = $cell.Value2
    if ($current -eq 46060) {
        $cell.Value2 = 46061
    }
}
